$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill the previously-empty row 11 with the new failed-spec remedy entry.
# (No rows are shifted - this just populates row 11 directly.)
$ws.Range("A11").Value = "SendingErrorInErrorVariableAndCallingWebservice"
$ws.Range("C11").Value = "Connection with credentials ""WarewolfAdmin"" and W@rEw0lf@dm1n fails on ""tst-ci-remote:3142"" server and returns ""(401) Unauthorized"" error"
$ws.Range("B11").Value = "Credentials on tst-ci-remote:3142 fails, need to correct the credentails (to be checked at tst-ci-remote:3142)"

$ws.Rows.Item(11).RowHeight = 28.8

# Update the selected cell to match the new active selection.
$ws.Range("A11").Select() | Out-Null
